$d = $word.ActiveDocument

# 1. Replace the text "OK." with "Asasasas".
$d.Content.Find.Execute("OK.", $false, $false, $false, $false, $false, $true, 1, $false, "Asasasas", 2)

# 2. Insert a brand new paragraph right after the "Asasasas" paragraph, carrying
#    the same paragraph/run formatting, and put "aassas" into it. The existing
#    _GoBack bookmark (which immediately followed the old "OK." run) stays put,
#    i.e. still attached to the "Asasasas" paragraph.
$found = $d.Range(0, 0)
$found.Find.Execute("Asasasas")
$found.Collapse(0)
$found.InsertParagraphAfter()

$again = $d.Range(0, 0)
$again.Find.Execute("Asasasas")
$para = $again.Paragraphs(1)
$newPara = $para.Next()
$newPara.Range.InsertAfter("aassas")
